$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values, regenerated from Strike# data (s_vals), per row (row 2..31)
$kValues = @{
    2  = 7
    3  = 5
    4  = 3
    5  = 3
    6  = 7
    7  = 4
    8  = 7
    9  = 7
    10 = 3
    11 = 5
    12 = 7
    13 = 8
    14 = 11
    15 = 6
    16 = 5
    17 = 4
    18 = 3
    19 = 3
    20 = 4
    21 = 2
    22 = 8
    23 = 7
    24 = 4
    25 = 6
    26 = 8
    27 = 0
    28 = 4
    29 = 10
    30 = 6
    31 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
